$d = $word.ActiveDocument

# The document repeated the same screenshot image in several standalone,
# highlighted (FFF3CD) paragraphs: once as the legitimate "featured"/intro
# capture, and again later (after "Presentación del proyecto" and after the
# "Pie de página" bullet) as accidental duplicates. Keep the first
# occurrence of each embedded image and remove every later paragraph that
# merely repeats it, deleting the whole paragraph (incl. its pPr/shading),
# exactly as the diff does.

$count = $d.InlineShapes.Count

# Resolve the embedded media filename for each inline picture so duplicates
# of the same image can be detected regardless of relationship-id reuse.
$mediaTargets = @()
for ($i = 1; $i -le $count; $i++) {
    $shape = $d.InlineShapes.Item($i)
    $xml = $shape.Range.WordOpenXML
    $target = ""
    $idx = $xml.IndexOf('Target="media')
    if ($idx -ge 0) {
        $segment = $xml.Substring($idx, 80)
        $startQuote = $segment.IndexOf('"') + 1
        $endQuote = $segment.IndexOf('"', $startQuote)
        $target = $segment.Substring($startQuote, $endQuote - $startQuote)
    }
    $mediaTargets += $target
}

# Walk from the end so deleting a duplicate paragraph never shifts the
# positions we still need to inspect.
for ($i = $count; $i -ge 1; $i--) {
    $seenBefore = @{}
    for ($j = 1; $j -lt $i; $j++) {
        $seenBefore[$mediaTargets[$j - 1]] = $true
    }

    $target = $mediaTargets[$i - 1]
    if ($target -ne "" -and $seenBefore.ContainsKey($target)) {
        $shape = $d.InlineShapes.Item($i)
        $shapeRange = $shape.Range
        $paraRange = $shapeRange.Paragraphs.Item(1).Range

        # Safety check: only delete if the picture is the paragraph's sole
        # content (the duplicated screenshots live alone in their own
        # highlighted paragraph).
        if ($paraRange.Start -eq $shapeRange.Start -and $paraRange.End -eq $shapeRange.End) {
            $paraRange.Delete()
        }
    }
}
